$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.428, -0.112, -0.2, -1, 38, -40),
    @(0.418, 0.259, 0.148, 165, -13, -154),
    @(-0.219, -0.516, 0.31, -38, 13, -9),
    @(0.162, 0.384, 0.107, 132, 28, -79),
    @(-0.296, 0.232, -0.135, 38, 34, -70),
    @(0.272, 0.235, -0.074, -82, 33, -154),
    @(-0.191, -0.194, 0.725, 65, 77, -121),
    @(0.179, 0.283, -0.161, -131, -51, -36),
    @(-0.242, 0.269, 0.167, 100, 49, 132),
    @(-0.154, -0.173, 0.489, -89, -5, -177)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}
